# Weekly Fruta/Hortaliza update: insert a new price record as the new
# most-recent row (row 92), pushing the existing rows 92-94 down to 93-95.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 92 (shifts old rows 92-94 -> 93-95).
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new record.
$ws.Cells.Item(92, 1).Value = 6
$ws.Cells.Item(92, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(92, 3).Value = "Metropolitana"
$ws.Cells.Item(92, 4).Value = 45239
$ws.Cells.Item(92, 5).Value = 13
$ws.Cells.Item(92, 6).Value = "Fruta"
$ws.Cells.Item(92, 7).Value = 100108
$ws.Cells.Item(92, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(92, 9).Value = 100108007
$ws.Cells.Item(92, 10).Value = "Coco"
$ws.Cells.Item(92, 11).Value = "Sin especificar"
$ws.Cells.Item(92, 12).Value = "Primera"
$ws.Cells.Item(92, 13).Value = 95
$ws.Cells.Item(92, 14).Value = 30000
$ws.Cells.Item(92, 15).Value = 30000
$ws.Cells.Item(92, 16).Value = 30000
$ws.Cells.Item(92, 17).Value = "`$/malla 20 unidades"
$ws.Cells.Item(92, 18).Value = "Perú"
$ws.Cells.Item(92, 19).Value = 1500
$ws.Cells.Item(92, 20).Value = 20
